$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7206.857
$ws.Range("J19").Value = 6487.75
$ws.Range("L19").Value = 6487.75
$ws.Range("N19").Value = -6837.75

$ws.Range("H64").Value = 10833.333
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 16666.666
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 16666.666
$ws.Range("M64").Value = -4752
$ws.Range("N64").Value = -17162.666

$ws.Range("H67").Value = 10833.333
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 16666.666
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 16666.666
$ws.Range("M67").Value = -4142
$ws.Range("N67").Value = -18382.666

$ws.Range("H96").Value = 1724.75
$ws.Range("I96").Value = 699.5
$ws.Range("J96").Value = 2750
$ws.Range("K96").Value = 2098.5
$ws.Range("L96").Value = 8250
$ws.Range("M96").Value = -725.5
$ws.Range("N96").Value = -10996

$ws.Range("H98").Value = 829.2632
$ws.Range("I98").Value = 826.1539
$ws.Range("K98").Value = 826.1539
$ws.Range("M98").Value = 671.8461

$ws.Range("H100").Value = 2109.1765
$ws.Range("I100").Value = 2459.5386
$ws.Range("J100").Value = 970.5
$ws.Range("K100").Value = 2459.5386
$ws.Range("L100").Value = 970.5
$ws.Range("M100").Value = -1918.5386
$ws.Range("N100").Value = -2052.5

$ws.Range("H106").Value = 1474
$ws.Range("I106").Value = 1474
$ws.Range("K106").Value = 1474
$ws.Range("M106").Value = -843

$ws.Range("H122").Value = 829.2632
$ws.Range("I122").Value = 826.1539
$ws.Range("K122").Value = 2478.4617
$ws.Range("M122").Value = -28.46169999999984

$ws.Range("H127").Value = 6830.3335
$ws.Range("I127").Value = 6830.3335
$ws.Range("K127").Value = 20491.0005
$ws.Range("M127").Value = -15531.0005

$ws.Range("H129").Value = 5081.5
$ws.Range("I129").Value = 5597.8
$ws.Range("J129").Value = 2500
$ws.Range("K129").Value = 16793.4
$ws.Range("L129").Value = 7500
$ws.Range("M129").Value = -11793.4
$ws.Range("N129").Value = -17500

$ws.Range("H132").Value = 4217.2964
$ws.Range("I132").Value = 3518.5293
$ws.Range("J132").Value = 5405.2
$ws.Range("K132").Value = 10555.5879
$ws.Range("L132").Value = 16215.6
$ws.Range("M132").Value = -8025.5879
$ws.Range("N132").Value = -21275.6

$ws.Range("H138").Value = 3428.7273
$ws.Range("I138").Value = 3079.6667
$ws.Range("K138").Value = 9239.000100000001
$ws.Range("M138").Value = -4099.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9505
$ws.Range("J2").Value = 10010
$ws.Range("L2").Value = 10010
$ws.Range("N2").Value = -10236

$ws.Range("H32").Value = 3304.5
$ws.Range("I32").Value = 2245.2
$ws.Range("J32").Value = 11249.25
$ws.Range("K32").Value = 2245.2
$ws.Range("L32").Value = 11249.25
$ws.Range("M32").Value = -1958.2
$ws.Range("N32").Value = -11823.25

$ws.Range("H74").Value = 1599.3636
$ws.Range("J74").Value = 1199
$ws.Range("L74").Value = 1199
$ws.Range("N74").Value = -2947

$ws.Range("H77").Value = 1599.3636
$ws.Range("J77").Value = 1199
$ws.Range("L77").Value = 5995
$ws.Range("N77").Value = -14731

$ws.Range("H97").Value = 926.95
$ws.Range("I97").Value = 784.9375
$ws.Range("J97").Value = 1495
$ws.Range("K97").Value = 784.9375
$ws.Range("L97").Value = 1495
$ws.Range("M97").Value = -288.9375
$ws.Range("N97").Value = -2487

$ws.Range("H110").Value = 1005.7
$ws.Range("I110").Value = 1005.7
$ws.Range("K110").Value = 1005.7
$ws.Range("M110").Value = 1039.3

$ws.Range("H116").Value = 9505
$ws.Range("J116").Value = 10010
$ws.Range("L116").Value = 10010
$ws.Range("N116").Value = -14598

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9505
$ws.Range("J3").Value = 10010
$ws.Range("L3").Value = 10010
$ws.Range("N3").Value = -10238

$ws.Range("H107").Value = 899.8
$ws.Range("I107").Value = 824.7143
$ws.Range("K107").Value = 824.7143
$ws.Range("M107").Value = 1095.2857

$ws.Range("H134").Value = 4526.375
$ws.Range("I134").Value = 5035.3335
$ws.Range("K134").Value = 15106.0005
$ws.Range("M134").Value = -12571.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2012.4445
$ws.Range("I31").Value = 1623
$ws.Range("K31").Value = 1623
$ws.Range("M31").Value = -1328

$ws.Range("H34").Value = 2012.4445
$ws.Range("I34").Value = 1623
$ws.Range("K34").Value = 1623
$ws.Range("M34").Value = -1421

$ws.Range("H105").Value = 2595
$ws.Range("I105").Value = 2595
$ws.Range("K105").Value = 2595
$ws.Range("M105").Value = -848

$ws.Range("H122").Value = 1838.6666
$ws.Range("I122").Value = 1838.6666
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5515.9998
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -3065.9998

$ws.Range("H132").Value = 4656.4165
$ws.Range("I132").Value = 5167.8
$ws.Range("K132").Value = 15503.4
$ws.Range("M132").Value = -12973.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 500000400
$ws.Range("I4").Value = 400000480
$ws.Range("J4").Value = 1000000000
$ws.Range("K4").Value = 1200001440
$ws.Range("L4").Value = 3000000000
$ws.Range("M4").Value = -1200001328
$ws.Range("N4").Value = -3000000224

$ws.Range("H68").Value = 409.5
$ws.Range("I68").Value = 320
$ws.Range("J68").Value = 499
$ws.Range("K68").Value = 960
$ws.Range("L68").Value = 1497
$ws.Range("N68").Value = -3119
$ws.Range("M68").Value = -149

$ws.Range("H71").Value = 409.5
$ws.Range("I71").Value = 320
$ws.Range("J71").Value = 499
$ws.Range("K71").Value = 2880
$ws.Range("L71").Value = 4491
$ws.Range("N71").Value = -12603
$ws.Range("M71").Value = 1176

$ws.Range("H109").Value = 1411.75
$ws.Range("I109").Value = 210.66667
$ws.Range("J109").Value = 5015
$ws.Range("K109").Value = 632.00001
$ws.Range("L109").Value = 15045
$ws.Range("M109").Value = 407.99999
$ws.Range("N109").Value = -17125

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

$ws.Range("H131").Value = 1137.25
$ws.Range("I131").Value = 1049.5
$ws.Range("J131").Value = 1166.5
$ws.Range("K131").Value = 3148.5
$ws.Range("L131").Value = 3499.5
$ws.Range("M131").Value = 1891.5
$ws.Range("N131").Value = -13579.5

$ws.Range("H137").Value = 3257.5
$ws.Range("I137").Value = 3257.5
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 9772.5
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -4672.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 54999.75
$ws.Range("I119").Value = 40000
$ws.Range("J119").Value = 99999
$ws.Range("K119").Value = 40000
$ws.Range("L119").Value = 99999
$ws.Range("N119").Value = -109675
$ws.Range("M119").Value = -35162

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1815.4615
$ws.Range("I46").Value = 1335.7142
$ws.Range("J46").Value = 2375.1667
$ws.Range("K46").Value = 1335.7142
$ws.Range("L46").Value = 2375.1667
$ws.Range("M46").Value = -1147.7142
$ws.Range("N46").Value = -2751.1667

$ws.Range("H68").Value = 2250.25
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 1
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 1
$ws.Range("M68").Value = -2251
$ws.Range("N68").Value = -1499

$ws.Range("H71").Value = 2250.25
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 1
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 5
$ws.Range("M71").Value = -11256
$ws.Range("N71").Value = -7493

$ws.Range("H132").Value = 3801.3157
$ws.Range("I132").Value = 3915.2
$ws.Range("K132").Value = 11745.6
$ws.Range("M132").Value = -9215.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3260.5454
$ws.Range("J81").Value = 3795.4
$ws.Range("L81").Value = 7590.8
$ws.Range("N81").Value = -9712.799999999999

$ws.Range("H84").Value = 3260.5454
$ws.Range("J84").Value = 3795.4
$ws.Range("L84").Value = 37954
$ws.Range("N84").Value = -48562

$ws.Range("H101").Value = 85897
$ws.Range("I101").Value = 85896
$ws.Range("K101").Value = 85896
$ws.Range("M101").Value = -82651

$ws.Range("H107").Value = 933.3333
$ws.Range("J107").Value = 800
$ws.Range("L107").Value = 2400
$ws.Range("N107").Value = -6240
